$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / recalculated means
$ws.Range("F3").Value = -9
$ws.Range("F7").Value = -3
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -12
$ws.Range("F11").Value = 7
$ws.Range("F14").Value = -3
$ws.Range("F16").Value = -12
$ws.Range("F18").Value = -5
